$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Add the new "2022-Q3" worksheet.
#
#    The cleanest way to get a new sheet that matches the look of the
#    existing quarterly sheets (sheetPr/outline flags, sheetFormatPr,
#    pageMargins, header style, column-A style, ...) is to duplicate the
#    most similar existing sheet ("2021-Q4", whose headers already read
#    基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名,
#    exactly like the new sheet needs) and then edit its data in place.
# ------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($q4)                       # duplicate placed immediately before "2021-Q4"
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q3"

# The new sheet only needs a single fund row, so drop the duplicated
# rows 3-7 that came along with the "2021-Q4" data.
$newSheet.Rows("3:7").Delete()

# Overwrite row 2 with the 2022-Q3 fund-holding detail. Values that look
# numeric but must stay textual (fund code / percentages, matching the
# source data) are entered with a leading apostrophe and then have their
# style reset to "Normal" so no stray numeric/text format sticks around.
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'005166"
$newSheet.Range("B2").Style = "Normal"
$newSheet.Range("C2").Value = "嘉实润和量化6个月定期开放混合"
$newSheet.Range("D2").Value = "'0.22"
$newSheet.Range("D2").Style = "Normal"
$newSheet.Range("E2").Value = "'24.64"
$newSheet.Range("E2").Style = "Normal"
$newSheet.Range("F2").Value = "'0.54"
$newSheet.Range("F2").Style = "Normal"
$newSheet.Range("G2").Value = "'0.0012"
$newSheet.Range("G2").Style = "Normal"
$newSheet.Range("H2").Value = 7

# ------------------------------------------------------------------
# 2) Update the "总计" (Total) summary sheet: insert a new row for the
#    2022-Q3 quarter above the existing 2021-Q4 row, shifting the rest
#    down by one row.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Insert a new blank row at row 2 (pushes 2021-Q4 / 2021-Q3 / 2021-Q2 down).
$total.Rows.Item(2).Insert()

# The freshly-inserted row auto-inherits a "painted" style across B:D from
# the row above/below; clear that back to the unstyled look the other data
# rows use, then copy the proper (bold/centered/bordered) style that column
# A data cells use onto the new A2 cell.
$total.Range("B2:D2").Style = "Normal"
$total.Range("A3").Copy($total.Range("A2"))

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0

# Renumber the running index in column A for the rows pushed down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# Restore "总计" as the active sheet (unchanged from the original file).
$total.Activate()
